# Updates cryptos list price/volume columns (D, E) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.523.07"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.466.79"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'314.28"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'91.59"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("D7").Value = "'0.545"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.512"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "'32.16"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "2.849.14"
$ws.Range("D14").Value = "'6.82"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "'16.03"
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").Value = "2.487.29"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "'0.766"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "41.487.12"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'6.48"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "0.0₃0946"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "'71.27"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'11.03"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "'235.75"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "'24.60"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'35.31"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'155.82"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("D32").Value = "'5.41"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'0.0755"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").Value = "'17.12"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("E36").Value = "  -8.20%  "
$ws.Range("D37").Value = "'2.86"
$ws.Range("E37").Value = "  -7.20%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("E41").Value = "  -5.06%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "1.950.69"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "'0.0282"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "'18.57"
$ws.Range("E45").Value = "  -4.68%  "
$ws.Range("D46").Value = "'2.91"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("D47").Value = "'9.04"
$ws.Range("E47").Value = "  +3.51%  "
$ws.Range("D48").Value = "2.707.22"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'96.91"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "'66.61"
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("E51").Value = "  -3.96%  "
